$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) columns for each coin row per the latest scrape.
# D-column values are forced to Text format so numeric-looking strings (e.g. "10.20",
# "0.00001015") keep their exact original text representation instead of being
# reinterpreted as numbers; the style is then reset to Normal so no visible formatting
# change is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.163.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.435.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9152"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "277.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3659"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.81%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3124"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.21"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.016"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06508"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.005"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.382"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.071"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.444.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001015"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9338"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05609"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.417"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.267"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.296.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.182"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "135.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.598.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.670"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8138"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.852"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07628"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.493"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +16.94%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05954"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.66%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.674"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.135"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.20"
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01985"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9312"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1822"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.898"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -18.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.524"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5223"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5125"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.761"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06323"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9985"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
